# Update countries & provincias Spain
# Applies the data refresh described in the commit:
#  - Updated "last updated" timestamp
#  - Kuwait's case counts overtook Sudafrica's -> rows swap order (text + numbers)
#  - Malta's case counts overtook Sierra Leona's -> rows swap order (text + numbers)
#  - Refreshed numeric totals for several other countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "datos actualizados" timestamp cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 13:35"

# 2. Estados Unidos (row 4)
$ws.Range("B4").Value = 1645646
$ws.Range("C4").Value = 552
$ws.Range("E4").Value = 1144755
$ws.Range("G4").Value = 16
$ws.Range("H4").Value = 97663

# 3. India (row 14)
$ws.Range("B14").Value = 126308
$ws.Range("C14").Value = 1514
$ws.Range("D14").Value = 52258
$ws.Range("E14").Value = 70296
$ws.Range("G14").Value = 28
$ws.Range("H14").Value = 3754

# 4. Kuwait overtakes Sudafrica - swap rows 37/38 (country names + data)
$ws.Range("A37").Value = "Kuwait"
$ws.Range("B37").Value = 20464
$ws.Range("C37").Value = 900
$ws.Range("D37").Value = 5747
$ws.Range("E37").Value = 14569
$ws.Range("G37").Value = 10
$ws.Range("H37").Value = 148

$ws.Range("A38").Value = "Sudafrica"
$ws.Range("B38").Value = 20125
$ws.Range("D38").Value = 10104
$ws.Range("E38").Value = 9624
$ws.Range("H38").Value = 397

# 5. Senegal (row 78)
$ws.Range("B78").Value = 2976
$ws.Range("C78").Value = 67
$ws.Range("D78").Value = 1416
$ws.Range("E78").Value = 1526
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 34

# 6. San Marino (row 125)
$ws.Range("B125").Value = 665
$ws.Range("C125").Value = 4
$ws.Range("D125").Value = 266
$ws.Range("E125").Value = 357
$ws.Range("G125").Value = 1
$ws.Range("H125").Value = 42

# 7. Malta overtakes Sierra Leona - swap rows 127/128 (country names + data)
$ws.Range("A127").Value = "Malta"
$ws.Range("B127").Value = 609
$ws.Range("C127").Value = 9
$ws.Range("D127").Value = 473
$ws.Range("E127").Value = 130
$ws.Range("H127").Value = 6

$ws.Range("A128").Value = "Sierra Leona"
$ws.Range("B128").Value = 606
$ws.Range("D128").Value = 230
$ws.Range("E128").Value = 338
$ws.Range("H128").Value = 38
